$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2933498322963715
$ws.Range("B1").Value = 1.938002228736877
$ws.Range("C1").Value = 4.34821605682373
$ws.Range("D1").Value = 1.599398612976074
$ws.Range("E1").Value = 1.048843264579773
